$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix A11 value (12 -> 13)
$ws.Range("A11").Value = 13

# New header for column J - copy formatting from I1 (bold/border/centered style)
# then set the text value.
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "Jena"

# New idf values for column J
$ws.Range("J2").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("J4").Value = 0.3333333333333333
$ws.Range("J5").Value = 0.5
$ws.Range("J6").Value = 0
$ws.Range("J7").Value = 0.5
$ws.Range("J8").Value = 0
$ws.Range("J9").Value = 1
$ws.Range("J10").Value = 1
$ws.Range("J11").Value = 0.5
